$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for paragraph $paraIndex : $findText"
    }
}

Replace-InParagraph 6 "Fornecer ao estudante os principais tipos de síntese orgânica e inorgânica de materiais bem como apresentar as principais técnicas analíticas para caracterização de materiais." "Introdução à química e sua associação com síntese de novos materiais. A visão moderna do átomo  e Ligações químicas. Estrutura cristalina e técnicas de caracterização cristalográfica. Filmes finos epitaxiais e filmes de uma maneira geral e seu impacto na tecnologica moderna. Crescimento de cristais  Materiais amorfos, síntese e aplicações. Processos e Técnicas de crescimento de cristais de um modo geral. Polímeros condutores e suas aplicações em tecnologica moderna."
Replace-InParagraph 7 "Provide the student with the main types of organic and inorganic synthesis of materials as well as presenting the main analytical techniques for material characterization." "Introduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology."
Replace-InParagraph 9 "5840730 - Antonio Jefferson da Silva Machado`v" "Fornecer ao estudante os principais tipos de síntese orgânica e inorgânica de materiais bem como apresentar as principais técnicas analíticas para caracterização de materiais.`v"
Replace-InParagraph 9 "5840897 - Clodoaldo Saron" "Química de materiais: definição; papel da química na ciência de materiais; fundamentos.`vAtomística e a visão moderna do átomo com fundamentos quânticos.Tipos de ligações químicas: forças de van der Waals, potencial de Lennard-Jones, ligação covalente, ligações por coordenação, ligações iônicas e ligações metálicas.`vMateriais policristalinos e monocristalinos. A ordem cristalográfica e técnicas de caracterização cristalográfica e microscópica. A importância de monocristais em aplicações eletrônicas. Técnicas de crescimento de cristais de alta qualidade tais como: método do fluxo, método Czochralski, método Brigdmann, método do transporte de vapor e método de crescimento de transporte de vapor modificado e isotérmico. Materiais amorfos e sua importância para a tecnologica moderna. Conceitos e técnicas de crescimento de materiais amorfos. Filmes finos epitaxiais, técnicas de crescimento tais como: vapor químico, sputtering, laser ablation e MBE. Filmes finos crescidos por eletrólise para revestimento protetivo, conceitos e aplicações. Síntese de polímeros condutores, conceitos e aplicações como dispositivos eletrônicos."
Replace-InParagraph 11 "Introdução à química e sua associação com síntese de novos materiais. A visão moderna do átomo  e Ligações químicas. Estrutura cristalina e técnicas de caracterização cristalográfica. Filmes finos epitaxiais e filmes de uma maneira geral e seu impacto na tecnologica moderna. Crescimento de cristais  Materiais amorfos, síntese e aplicações. Processos e Técnicas de crescimento de cristais de um modo geral. Polímeros condutores e suas aplicações em tecnologica moderna." "Aulas expositivas e práticas ministradas em laboratório."
Replace-InParagraph 12 "Introduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology." "Provide the student with the main types of organic and inorganic synthesis of materials as well as presenting the main analytical techniques for material characterization."
Replace-InParagraph 14 "Química de materiais: definição; papel da química na ciência de materiais; fundamentos.`vAtomística e a visão moderna do átomo com fundamentos quânticos.Tipos de ligações químicas: forças de van der Waals, potencial de Lennard-Jones, ligação covalente, ligações por coordenação, ligações iônicas e ligações metálicas.`vMateriais policristalinos e monocristalinos. A ordem cristalográfica e técnicas de caracterização cristalográfica e microscópica. A importância de monocristais em aplicações eletrônicas. Técnicas de crescimento de cristais de alta qualidade tais como: método do fluxo, método Czochralski, método Brigdmann, método do transporte de vapor e método de crescimento de transporte de vapor modificado e isotérmico. Materiais amorfos e sua importância para a tecnologica moderna. Conceitos e técnicas de crescimento de materiais amorfos. Filmes finos epitaxiais, técnicas de crescimento tais como: vapor químico, sputtering, laser ablation e MBE. Filmes finos crescidos por eletrólise para revestimento protetivo, conceitos e aplicações. Síntese de polímeros condutores, conceitos e aplicações como dispositivos eletrônicos." "Média simples de duas provas escritas,  Conceito Final = (P1 + P2)/2"
Replace-InParagraph 17 "Média simples de duas provas escritas,  Conceito Final = (P1 + P2)/2`v" "ALLCOCK, H. R. Introduction to Materials Chemistry, Wiley, Nova Iorque, 2008.`vFAHLMAN, B. D. Materials Chemistry, Springer, Holanda, 2007.`vZHANG, S.; LI, L.; KUMAR, A. Materials Characterization Techniques, Boca Raton: CRC Press, 2008.`vLENG, Y. Materials Characterization: Introduction to Microscopic and Spectroscopic Methods, Wiley, Cingapura, 2008.`v"
Replace-InParagraph 17 "Aplicação de duas provas escritas dentro do prazo regimental antes do início do próximo semestre letivo." "5840730 - Antonio Jefferson da Silva Machado"
Replace-InParagraph 17 "Aulas expositivas e práticas ministradas em laboratório.`v" "Aplicação de duas provas escritas dentro do prazo regimental antes do início do próximo semestre letivo.`v"
Replace-InParagraph 19 "ALLCOCK, H. R. Introduction to Materials Chemistry, Wiley, Nova Iorque, 2008.`vFAHLMAN, B. D. Materials Chemistry, Springer, Holanda, 2007.`vZHANG, S.; LI, L.; KUMAR, A. Materials Characterization Techniques, Boca Raton: CRC Press, 2008.`vLENG, Y. Materials Characterization: Introduction to Microscopic and Spectroscopic Methods, Wiley, Cingapura, 2008." "5840897 - Clodoaldo Saron"
